$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-29 Tuesday" "2025-04-30 Wednesday"

Replace-Text "768÷2=384, 0" "510÷8=63, 6"
Replace-Text "153÷5=30, 3" "103÷4=25, 3"
Replace-Text "207÷5=41, 2" "737÷2=368, 1"
Replace-Text "421÷8=52, 5" "171÷2=85, 1"
Replace-Text "764÷7=109, 1" "403÷2=201, 1"

Replace-Text "221÷4=55, 1" "478÷8=59, 6"
Replace-Text "558÷9=62, 0" "334÷5=66, 4"
Replace-Text "268÷6=44, 4" "100÷9=11, 1"
Replace-Text "806÷8=100, 6" "286÷8=35, 6"
Replace-Text "399÷5=79, 4" "838÷6=139, 4"

Replace-Text "868÷7=124, 0" "925÷7=132, 1"
Replace-Text "850÷9=94, 4" "146÷3=48, 2"
Replace-Text "707÷7=101, 0" "231÷7=33, 0"
Replace-Text "372÷5=74, 2" "275÷8=34, 3"
Replace-Text "688÷5=137, 3" "758÷2=379, 0"

Replace-Text "423÷3=141, 0" "372÷2=186, 0"
Replace-Text "579÷8=72, 3" "436÷5=87, 1"
Replace-Text "514÷5=102, 4" "580÷4=145, 0"
Replace-Text "852÷3=284, 0" "928÷8=116, 0"
Replace-Text "850÷2=425, 0" "894÷5=178, 4"

Replace-Text "698÷7=99, 5" "240÷2=120, 0"
Replace-Text "664÷6=110, 4" "279÷3=93, 0"
Replace-Text "171÷3=57, 0" "352÷2=176, 0"
Replace-Text "304÷8=38, 0" "192÷6=32, 0"
Replace-Text "575÷9=63, 8" "224÷3=74, 2"
